$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F4").Value = 216
    $ws.Range("F5").Value = 2641
    $ws.Range("F7").Value = 363
}

$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F9").Value = 933

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F10").Value = 933
